$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 56
$ws.Range("A4").Value = 79
$ws.Range("A5").Value = 0
$ws.Range("A13").Value = 1
$ws.Range("A15").Value = 0
$ws.Range("A17").Value = 3
$ws.Range("A18").Value = 2
